# Update microstate list: remove replicate microstate (v1.4.1)
#
# The sheet lists one microstate (ID + canonical isomeric SMILES) per row,
# each row carrying an associated 2D-depiction picture anchored to it.
# Row 28 (microstate "SM10_micro029") is a replicate entry and must be
# removed, together with its 2D-depiction picture. Deleting the worksheet
# row shifts all subsequent rows up by one (and the now-unreferenced
# shared strings for the removed ID/SMILES drop out of the shared string
# table), while the picture anchored at the very bottom of the drawing
# (the one accompanying the last row, "Picture 32") is deleted to keep
# the pictures aligned one-per-row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the replicate microstate row ("SM10_micro029") entirely; this
# shifts rows 29-34 up to become rows 28-33.
$ws.Rows.Item(28).Delete()

# Remove the now-surplus 2D-depiction picture (the last one in the
# drawing, originally paired with the last data row).
$ws.Shapes.Item("Picture 32").Delete()
